$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.80670055787346
$ws.Range("D2").Value = 8.328752107955303
$ws.Range("E2").Value = 13.87563784222323
$ws.Range("F2").Value = 36.15111055731585
$ws.Range("G2").Value = 43.14642794347948
$ws.Range("H2").Value = 17.46453111729555
$ws.Range("I2").Value = 25.37782523510359
$ws.Range("J2").Value = 10.29997905330501
$ws.Range("L2").Value = 14.47518540856798
$ws.Range("B3").Value = 19.25402173270246
$ws.Range("D3").Value = 8.233686616986342
$ws.Range("E3").Value = 13.6637665262217
$ws.Range("F3").Value = 36.24921599958065
$ws.Range("G3").Value = 42.95026760909532
$ws.Range("H3").Value = 17.51066600940263
$ws.Range("I3").Value = 25.61999666579085
$ws.Range("J3").Value = 10.2267740057106
$ws.Range("L3").Value = 14.14236088811965
$ws.Range("B4").Value = 18.90736832221653
$ws.Range("D4").Value = 8.174165538885109
$ws.Range("E4").Value = 13.53163829552898
$ws.Range("F4").Value = 36.32374210454149
$ws.Range("G4").Value = 42.85127095485308
$ws.Range("H4").Value = 17.54447411407039
$ws.Range("I4").Value = 25.77611501325435
$ws.Range("J4").Value = 10.1818199267821
$ws.Range("L4").Value = 13.93516014864188
$ws.Range("B5").Value = 18.76445215284522
$ws.Range("D5").Value = 8.149632596958224
$ws.Range("E5").Value = 13.47732617077067
$ws.Range("F5").Value = 36.35767788893053
$ws.Range("G5").Value = 42.81633713851008
$ws.Range("H5").Value = 17.55962144269306
$ws.Range("I5").Value = 25.8416066107342
$ws.Range("J5").Value = 10.16350941759291
$ws.Range("L5").Value = 13.85011579730696
$ws.Range("B6").Value = 18.74062690634602
$ws.Range("D6").Value = 8.145542497850361
$ws.Range("E6").Value = 13.46828066939792
$ws.Range("F6").Value = 36.36352743341933
$ws.Range("G6").Value = 42.81086326288181
$ws.Range("H6").Value = 17.5622191620169
$ws.Range("I6").Value = 25.85259466668624
$ws.Range("J6").Value = 10.16046978746984
$ws.Range("L6").Value = 13.83596072458522
$ws.Range("B7").Value = 18.90544734023119
$ws.Range("D7").Value = 8.17383578751754
$ws.Range("E7").Value = 13.53090766271011
$ws.Range("F7").Value = 36.32418537422082
$ws.Range("G7").Value = 42.85077791651679
$ws.Range("H7").Value = 17.54467285918077
$ws.Range("I7").Value = 25.7769906688722
$ws.Range("J7").Value = 10.18157293654322
$ws.Range("L7").Value = 13.93401552913011
$ws.Range("B8").Value = 19.61776049547402
$ws.Range("D8").Value = 8.296219250381322
$ws.Range("E8").Value = 13.80303358955159
$ws.Range("F8").Value = 36.18195576442022
$ws.Range("G8").Value = 43.0743532301158
$ws.Range("H8").Value = 17.47929613844778
$ws.Range("I8").Value = 25.4597882223662
$ws.Range("J8").Value = 10.27474526659045
$ws.Range("L8").Value = 14.36108075040539
$ws.Range("B9").Value = 20.94915034238902
$ws.Range("D9").Value = 8.526584753702988
$ws.Range("E9").Value = 14.31863798680183
$ws.Range("F9").Value = 36.01755447005701
$ws.Range("G9").Value = 43.68173047359871
$ws.Range("H9").Value = 17.39492120286884
$ws.Range("I9").Value = 24.89642343406171
$ws.Range("J9").Value = 10.45696391764966
$ws.Range("L9").Value = 15.17153329366813
$ws.Range("B10").Value = 21.87833224092905
$ws.Range("D10").Value = 8.689317179082183
$ws.Range("E10").Value = 14.6840679786382
$ws.Range("F10").Value = 35.96799503580236
$ws.Range("G10").Value = 44.22855115381777
$ws.Range("H10").Value = 17.36008275291013
$ws.Range("I10").Value = 24.51795258428752
$ws.Range("J10").Value = 10.58994433134802
$ws.Range("L10").Value = 15.74494072079803
$ws.Range("B11").Value = 22.28874488144628
$ws.Range("D11").Value = 8.761801646931159
$ws.Range("E11").Value = 14.84692962607033
$ws.Range("F11").Value = 35.96117317059532
$ws.Range("G11").Value = 44.49841692029318
$ws.Range("H11").Value = 17.35021125791974
$ws.Range("I11").Value = 24.35340134705515
$ws.Range("J11").Value = 10.65011946280194
$ws.Range("L11").Value = 15.9999577249499
$ws.Range("B12").Value = 22.44227086789121
$ws.Range("D12").Value = 8.789016978043193
$ws.Range("E12").Value = 14.90807830318919
$ws.Range("F12").Value = 35.96086826038925
$ws.Range("G12").Value = 44.60356786203717
$ws.Range("H12").Value = 17.3473383684502
$ws.Range("I12").Value = 24.29218043609276
$ws.Range("J12").Value = 10.67284975770423
$ws.Range("L12").Value = 16.09560803835678
$ws.Range("B13").Value = 22.40929211444166
$ws.Range("D13").Value = 8.783166202457235
$ws.Range("E13").Value = 14.89493271363633
$ws.Range("F13").Value = 35.96083234594003
$ws.Range("G13").Value = 44.58079146729752
$ws.Range("H13").Value = 17.34791852962071
$ws.Range("I13").Value = 24.3053169950667
$ws.Range("J13").Value = 10.66795707802788
$ws.Range("L13").Value = 16.07505013024049
$ws.Range("B14").Value = 22.30141406365669
$ws.Range("D14").Value = 8.764045399433657
$ws.Range("E14").Value = 14.85197106874809
$ws.Range("F14").Value = 35.96110233243178
$ws.Range("G14").Value = 44.50700893161456
$ws.Range("H14").Value = 17.34995753210603
$ws.Range("I14").Value = 24.34834282301012
$ws.Range("J14").Value = 10.65199069700236
$ws.Range("L14").Value = 16.00784579484639
$ws.Range("B15").Value = 22.23508627925907
$ws.Range("D15").Value = 8.75230269550786
$ws.Range("E15").Value = 14.82558654218837
$ws.Range("F15").Value = 35.9615648937988
$ws.Range("G15").Value = 44.46219784697478
$ws.Range("H15").Value = 17.3513193198754
$ws.Range("I15").Value = 24.37483935096059
$ws.Range("J15").Value = 10.64220311745091
$ws.Range("L15").Value = 15.9665592122251
$ws.Range("B16").Value = 21.85125297518913
$ws.Range("D16").Value = 8.684548141308776
$ws.Range("E16").Value = 14.67335337345803
$ws.Range("F16").Value = 35.96875879532205
$ws.Range("G16").Value = 44.2113332368878
$ws.Range("H16").Value = 17.36084865376779
$ws.Range("I16").Value = 24.52885930243973
$ws.Range("J16").Value = 10.58600440178716
$ws.Range("L16").Value = 15.72815040809541
$ws.Range("B17").Value = 21.61254786241359
$ws.Range("D17").Value = 8.642579667336735
$ws.Range("E17").Value = 14.57907120058311
$ws.Range("F17").Value = 35.97721225457462
$ws.Range("G17").Value = 44.062791465872
$ws.Range("H17").Value = 17.36822968111278
$ws.Range("I17").Value = 24.62529339261371
$ws.Range("J17").Value = 10.55143896497616
$ws.Range("L17").Value = 15.58034184730139
$ws.Range("B18").Value = 21.47410199032947
$ws.Range("D18").Value = 8.618295883405059
$ws.Range("E18").Value = 14.52452715827184
$ws.Range("F18").Value = 35.98355349519048
$ws.Range("G18").Value = 43.97934678609568
$ws.Range("H18").Value = 17.37303729958404
$ws.Range("I18").Value = 24.68147692905075
$ws.Range("J18").Value = 10.53152852810616
$ws.Range("L18").Value = 15.4947820646509
$ws.Range("B19").Value = 21.42703320289652
$ws.Range("D19").Value = 8.610049313672699
$ws.Range("E19").Value = 14.50600650387269
$ws.Range("F19").Value = 35.98595398316312
$ws.Range("G19").Value = 43.951438302008
$ws.Range("H19").Value = 17.37476145635888
$ws.Range("I19").Value = 24.70062302040625
$ws.Range("J19").Value = 10.52478249685537
$ws.Range("L19").Value = 15.46572209221202
$ws.Range("B20").Value = 21.6380782764924
$ws.Range("D20").Value = 8.647062332546311
$ws.Range("E20").Value = 14.58914062780221
$ws.Range("F20").Value = 35.97615916066124
$ws.Range("G20").Value = 44.07839829636864
$ws.Range("H20").Value = 17.36738572384641
$ws.Range("I20").Value = 24.61495362306536
$ws.Range("J20").Value = 10.55512162796151
$ws.Range("L20").Value = 15.59613325035466
$ws.Range("B21").Value = 22.33315263173616
$ws.Range("D21").Value = 8.769668051462991
$ws.Range("E21").Value = 14.86460444853451
$ws.Range("F21").Value = 35.96096106982362
$ws.Range("G21").Value = 44.52860101604628
$ws.Range("H21").Value = 17.34933510245897
$ws.Range("I21").Value = 24.33567551616874
$ws.Range("J21").Value = 10.65668203633306
$ws.Range("L21").Value = 16.02761089471119
$ws.Range("B22").Value = 22.77636844259575
$ws.Range("D22").Value = 8.848435231453228
$ws.Range("E22").Value = 15.04156822715125
$ws.Range("F22").Value = 35.96431467047503
$ws.Range("G22").Value = 44.84003772268202
$ws.Range("H22").Value = 17.34258344917662
$ws.Range("I22").Value = 24.15950885947774
$ws.Range("J22").Value = 10.72272201083838
$ws.Range("L22").Value = 16.30421619219463
$ws.Range("B23").Value = 22.54086490600899
$ws.Range("D23").Value = 8.806523945660826
$ws.Range("E23").Value = 14.94741218828722
$ws.Range("F23").Value = 35.96130398043607
$ws.Range("G23").Value = 44.6722718441934
$ws.Range("H23").Value = 17.34572352195971
$ws.Range("I23").Value = 24.25295195316183
$ws.Range("J23").Value = 10.68750953471025
$ws.Range("L23").Value = 16.15710506053298
$ws.Range("B24").Value = 21.62653973920443
$ws.Range("D24").Value = 8.645036203265381
$ws.Range("E24").Value = 14.5845892957489
$ws.Range("F24").Value = 35.9766306520576
$ws.Range("G24").Value = 44.07133635636304
$ws.Range("H24").Value = 17.36776552011695
$ws.Range("I24").Value = 24.61962592047674
$ws.Range("J24").Value = 10.553456814512
$ws.Range("L24").Value = 15.58899576572218
$ws.Range("B25").Value = 20.59697111672174
$ws.Range("D25").Value = 8.465363425814937
$ws.Range("E25").Value = 14.18133662893487
$ws.Range("F25").Value = 36.04961553202394
$ws.Range("G25").Value = 43.4995383447858
$ws.Range("H25").Value = 17.41300796421039
$ws.Range("I25").Value = 25.04258253443787
$ws.Range("J25").Value = 10.40778303308203
$ws.Range("L25").Value = 14.9557740931688
